$d = $word.ActiveDocument

# The memo table is the only table in the document; grab it and add a new
# row that inherits the formatting (incl. bullet numbering) of the last row.
$t = $d.Tables.Item(1)
$t.Rows.Add() | Out-Null

$newIndex = $t.Rows.Count
$dateCell = $t.Cell($newIndex, 1)
$taskCell = $t.Cell($newIndex, 2)

$dateCell.Range.Text = "18/08/2020"
$taskCell.Range.Text = "Mejora prototipo ciudad"
